# Insert a new weekly price record at row 354 of the "Pepino ensalada" sheet.
# Existing rows 354-368 shift down to 355-369 (unchanged), and the freshly
# inserted row 354 is populated with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 354..368 down to 355..369, leaving a blank row 354 to fill in.
$ws.Rows.Item(354).Insert()

$ws.Range("A354").Value = 4
$ws.Range("B354").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C354").Value = "Los Lagos"
$ws.Range("D354").Value = 44939
$ws.Range("E354").Value = 10
$ws.Range("F354").Value = 100112043
$ws.Range("G354").Value = "Pepino ensalada"
$ws.Range("H354").Value = "Sin especificar"
$ws.Range("I354").Value = "Primera"
$ws.Range("J354").Value = 400
$ws.Range("K354").Value = 18000
$ws.Range("L354").Value = 20000
$ws.Range("M354").Value = 19000
$ws.Range("N354").Value = "`$/caja 60 unidades"
$ws.Range("O354").Value = "Región de Arica y Parinacota"
$ws.Range("P354").Value = 317
$ws.Range("Q354").Value = 60
$ws.Range("R354").Value = "Hortaliza"
